$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "SubSeason" column (D) with per-row values derived from Season,
# followed by the header - matches shared-string insertion order
# (Active, Hibernation, SubSeason)
$ws.Range("D2").Value = "Active"
$ws.Range("D3").Value = "Active"
$ws.Range("D4").Value = "Hibernation"
$ws.Range("D5").Value = "Hibernation"
$ws.Range("D6").Value = "Active"
$ws.Range("D7").Value = "Active"
$ws.Range("D8").Value = "Hibernation"
$ws.Range("D9").Value = "Active"
$ws.Range("D10").Value = "Active"
$ws.Range("D1").Value = "SubSeason"

# Update selection to reflect the newly added column
$ws.Range("D1").Select()
